# "Ajustes de datos de biomasa" - Cambios en recurso aprovechable
# Update the "Regimen de aprovechamiento" (column I, "Recomendado") values
# on the "Datos" sheet for bovinos, porcinos, aves, ovinos and caña rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

$ws.Range("I3").Value = 1500
$ws.Range("I5").Value = 5000
$ws.Range("I7").Value = 100000
$ws.Range("I13").Value = 50
$ws.Range("I16").Value = 10

# Make "Datos" the active sheet/tab and leave the selection on the last
# edited cell (I7), matching the workbook's saved view state.
$ws.Activate()
$ws.Range("I7").Select()
